$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell F1 with value and copy header style from B1
$ws.Range("F1").Value = "time_taken"
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate F2:F117 with time_taken values
$arr = New-Object "object[,]" 116,1
$arr[0,0] = "2021-10-05 10:52:55.172760"
$arr[1,0] = "2021-10-05 10:52:55.172772"
$arr[2,0] = "2021-10-05 10:52:55.172775"
$arr[3,0] = "2021-10-05 10:52:55.172778"
$arr[4,0] = "2021-10-05 10:52:55.172781"
$arr[5,0] = "2021-10-05 10:52:55.172784"
$arr[6,0] = "2021-10-05 10:52:55.172786"
$arr[7,0] = "2021-10-05 10:52:55.172789"
$arr[8,0] = "2021-10-05 10:52:55.172791"
$arr[9,0] = "2021-10-05 10:52:55.172794"
$arr[10,0] = "2021-10-05 10:52:55.172796"
$arr[11,0] = "2021-10-05 10:52:55.172799"
$arr[12,0] = "2021-10-05 10:52:55.172801"
$arr[13,0] = "2021-10-05 10:52:55.172804"
$arr[14,0] = "2021-10-05 10:52:55.172806"
$arr[15,0] = "2021-10-05 10:52:55.172809"
$arr[16,0] = "2021-10-05 10:52:55.172811"
$arr[17,0] = "2021-10-05 10:52:55.172814"
$arr[18,0] = "2021-10-05 10:52:55.172816"
$arr[19,0] = "2021-10-05 10:52:55.172819"
$arr[20,0] = "2021-10-05 10:52:55.172821"
$arr[21,0] = "2021-10-05 10:52:55.172824"
$arr[22,0] = "2021-10-05 10:52:55.172826"
$arr[23,0] = "2021-10-05 10:52:55.172829"
$arr[24,0] = "2021-10-05 10:52:55.172831"
$arr[25,0] = "2021-10-05 10:52:55.172834"
$arr[26,0] = "2021-10-05 10:52:55.172837"
$arr[27,0] = "2021-10-05 10:52:55.172839"
$arr[28,0] = "2021-10-05 10:52:55.172841"
$arr[29,0] = "2021-10-05 10:52:55.172844"
$arr[30,0] = "2021-10-05 10:52:55.172846"
$arr[31,0] = "2021-10-05 10:52:55.172849"
$arr[32,0] = "2021-10-05 10:52:55.172852"
$arr[33,0] = "2021-10-05 10:52:55.172854"
$arr[34,0] = "2021-10-05 10:52:55.172857"
$arr[35,0] = "2021-10-05 10:52:55.172860"
$arr[36,0] = "2021-10-05 10:52:55.172862"
$arr[37,0] = "2021-10-05 10:52:55.172865"
$arr[38,0] = "2021-10-05 10:52:55.172867"
$arr[39,0] = "2021-10-05 10:52:55.172870"
$arr[40,0] = "2021-10-05 10:52:55.172873"
$arr[41,0] = "2021-10-05 10:52:55.172875"
$arr[42,0] = "2021-10-05 10:52:55.172878"
$arr[43,0] = "2021-10-05 10:52:55.172880"
$arr[44,0] = "2021-10-05 10:52:55.172883"
$arr[45,0] = "2021-10-05 10:52:55.172885"
$arr[46,0] = "2021-10-05 10:52:55.172888"
$arr[47,0] = "2021-10-05 10:52:55.172890"
$arr[48,0] = "2021-10-05 10:52:55.172893"
$arr[49,0] = "2021-10-05 10:52:55.172895"
$arr[50,0] = "2021-10-05 10:52:55.172898"
$arr[51,0] = "2021-10-05 10:52:55.172900"
$arr[52,0] = "2021-10-05 10:52:55.172903"
$arr[53,0] = "2021-10-05 10:52:55.172905"
$arr[54,0] = "2021-10-05 10:52:55.172908"
$arr[55,0] = "2021-10-05 10:52:55.172910"
$arr[56,0] = "2021-10-05 10:52:55.172913"
$arr[57,0] = "2021-10-05 10:52:55.172915"
$arr[58,0] = "2021-10-05 10:52:55.172918"
$arr[59,0] = "2021-10-05 10:52:55.172920"
$arr[60,0] = "2021-10-05 10:52:55.172923"
$arr[61,0] = "2021-10-05 10:52:55.172925"
$arr[62,0] = "2021-10-05 10:52:55.172927"
$arr[63,0] = "2021-10-05 10:52:55.172930"
$arr[64,0] = "2021-10-05 10:52:55.172933"
$arr[65,0] = "2021-10-05 10:52:55.172936"
$arr[66,0] = "2021-10-05 10:52:55.172939"
$arr[67,0] = "2021-10-05 10:52:55.172941"
$arr[68,0] = "2021-10-05 10:52:55.172944"
$arr[69,0] = "2021-10-05 10:52:55.172946"
$arr[70,0] = "2021-10-05 10:52:55.172954"
$arr[71,0] = "2021-10-05 10:52:55.172957"
$arr[72,0] = "2021-10-05 10:52:55.172960"
$arr[73,0] = "2021-10-05 10:52:55.172962"
$arr[74,0] = "2021-10-05 10:52:55.172965"
$arr[75,0] = "2021-10-05 10:52:55.172967"
$arr[76,0] = "2021-10-05 10:52:55.172971"
$arr[77,0] = "2021-10-05 10:52:55.172974"
$arr[78,0] = "2021-10-05 10:52:55.172977"
$arr[79,0] = "2021-10-05 10:52:55.172979"
$arr[80,0] = "2021-10-05 10:52:55.172982"
$arr[81,0] = "2021-10-05 10:52:55.172984"
$arr[82,0] = "2021-10-05 10:52:55.172987"
$arr[83,0] = "2021-10-05 10:52:55.172989"
$arr[84,0] = "2021-10-05 10:52:55.172992"
$arr[85,0] = "2021-10-05 10:52:55.172994"
$arr[86,0] = "2021-10-05 10:52:55.172997"
$arr[87,0] = "2021-10-05 10:52:55.172999"
$arr[88,0] = "2021-10-05 10:52:55.173002"
$arr[89,0] = "2021-10-05 10:52:55.173004"
$arr[90,0] = "2021-10-05 10:52:55.173007"
$arr[91,0] = "2021-10-05 10:52:55.173009"
$arr[92,0] = "2021-10-05 10:52:55.173013"
$arr[93,0] = "2021-10-05 10:52:55.173016"
$arr[94,0] = "2021-10-05 10:52:55.173018"
$arr[95,0] = "2021-10-05 10:52:55.173021"
$arr[96,0] = "2021-10-05 10:52:55.173023"
$arr[97,0] = "2021-10-05 10:52:55.173026"
$arr[98,0] = "2021-10-05 10:52:55.173029"
$arr[99,0] = "2021-10-05 10:52:55.173033"
$arr[100,0] = "2021-10-05 10:52:55.173037"
$arr[101,0] = "2021-10-05 10:52:55.173041"
$arr[102,0] = "2021-10-05 10:52:55.173046"
$arr[103,0] = "2021-10-05 10:52:55.173050"
$arr[104,0] = "2021-10-05 10:52:55.173055"
$arr[105,0] = "2021-10-05 10:52:55.173059"
$arr[106,0] = "2021-10-05 10:52:55.173063"
$arr[107,0] = "2021-10-05 10:52:55.173067"
$arr[108,0] = "2021-10-05 10:52:55.173074"
$arr[109,0] = "2021-10-05 10:52:55.173080"
$arr[110,0] = "2021-10-05 10:52:55.173084"
$arr[111,0] = "2021-10-05 10:52:55.173088"
$arr[112,0] = "2021-10-05 10:52:55.173092"
$arr[113,0] = "2021-10-05 10:52:55.173097"
$arr[114,0] = "2021-10-05 10:52:55.173100"
$arr[115,0] = "2021-10-05 10:52:55.173102"
$ws.Range("F2:F117").Value = $arr

# Update the used dimension reference
$ws.UsedRange | Out-Null
